$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rushing" (sheet1): update a few players' rushing-attempt stats from
# the newly simulated Wild Card round.
# ---------------------------------------------------------------------------
$wsRushing = $wb.Worksheets.Item("Rushing")

# Row 2 - R.Tannehill
$wsRushing.Range("D2").Value = 7

# Row 6 - D.Foreman
$wsRushing.Range("C6").Value = 80
$wsRushing.Range("D6").Value = 52
$wsRushing.Range("F6").Value = 27

# Row 7 - D.Hilliard
$wsRushing.Range("C7").Value = 35
$wsRushing.Range("D7").Value = 28
$wsRushing.Range("E7").Value = 13
$wsRushing.Range("F7").Value = 11

# ---------------------------------------------------------------------------
# Sheet "Receiving" (sheet2): a new player (J.Jones) who appeared in the Wild
# Card round is logged, pushing every player below him down by one row, and
# several players' receiving stats are updated with the round's results.
# ---------------------------------------------------------------------------
$wsReceiving = $wb.Worksheets.Item("Receiving")

# Row 3 - K.Blasingame
$wsReceiving.Range("C3").Value = 2
$wsReceiving.Range("D3").Value = 2

# Row 5 - D.Foreman
$wsReceiving.Range("C5").Value = 10
$wsReceiving.Range("D5").Value = 8
$wsReceiving.Range("G5").Value = 2
$wsReceiving.Range("H5").Value = 2

# Row 6 - D.Hilliard
$wsReceiving.Range("C6").Value = 25
$wsReceiving.Range("D6").Value = 18

# Row 7 - A.Brown
$wsReceiving.Range("C7").Value = 101
$wsReceiving.Range("D7").Value = 74
$wsReceiving.Range("E7").Value = 25
$wsReceiving.Range("F7").Value = 18
$wsReceiving.Range("G7").Value = 12
$wsReceiving.Range("H7").Value = 10

# Give the new row 18 the same look (borders/alignment) as the existing
# table rows before filling in its values below.
$wsReceiving.Range("A17").Copy()
$wsReceiving.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 8 - new player J.Jones, logged from the Wild Card round
$wsReceiving.Range("A8").Value = 6
$wsReceiving.Range("B8").Value = "J.Jones"
$wsReceiving.Range("C8").Value = 6
$wsReceiving.Range("D8").Value = 5
$wsReceiving.Range("E8").Value = 3
$wsReceiving.Range("F8").Value = 0
$wsReceiving.Range("G8").Value = 1
$wsReceiving.Range("H8").Value = 1

# Row 9 - C.Rogers (shifted down from row 8, stats updated)
$wsReceiving.Range("A9").Value = 7
$wsReceiving.Range("B9").Value = "C.Rogers"
$wsReceiving.Range("C9").Value = 31
$wsReceiving.Range("D9").Value = 22
$wsReceiving.Range("E9").Value = 4
$wsReceiving.Range("F9").Value = 2
$wsReceiving.Range("G9").Value = 3
$wsReceiving.Range("H9").Value = 1

# Row 10 - R.McMath (shifted down from row 9, stats updated)
$wsReceiving.Range("A10").Value = 8
$wsReceiving.Range("B10").Value = "R.McMath"
$wsReceiving.Range("C10").Value = 3
$wsReceiving.Range("D10").Value = 2
$wsReceiving.Range("E10").Value = 1
$wsReceiving.Range("F10").Value = 0
$wsReceiving.Range("G10").Value = 1
$wsReceiving.Range("H10").Value = 1

# Row 11 - M.Johnson (shifted down from row 10, stats updated)
$wsReceiving.Range("A11").Value = 9
$wsReceiving.Range("B11").Value = "M.Johnson"
$wsReceiving.Range("C11").Value = 13
$wsReceiving.Range("D11").Value = 6
$wsReceiving.Range("E11").Value = 6
$wsReceiving.Range("F11").Value = 3
$wsReceiving.Range("G11").Value = 3
$wsReceiving.Range("H11").Value = 1

# Row 12 - D.Fitzpatrick (shifted down from row 11, unchanged stats)
$wsReceiving.Range("A12").Value = 10
$wsReceiving.Range("B12").Value = "D.Fitzpatrick"
$wsReceiving.Range("C12").Value = 6
$wsReceiving.Range("D12").Value = 5
$wsReceiving.Range("E12").Value = 2
$wsReceiving.Range("F12").Value = 0
$wsReceiving.Range("G12").Value = 2
$wsReceiving.Range("H12").Value = 2

# Row 13 - N.Westbrook-Ikhine (shifted down from row 12, stats updated)
$wsReceiving.Range("A13").Value = 11
$wsReceiving.Range("B13").Value = "N.Westbrook-Ikhine"
$wsReceiving.Range("C13").Value = 27
$wsReceiving.Range("D13").Value = 21
$wsReceiving.Range("E13").Value = 11
$wsReceiving.Range("F13").Value = 5
$wsReceiving.Range("G13").Value = 5
$wsReceiving.Range("H13").Value = 5

# Row 14 - C.Hollister (shifted down from row 13, stats updated)
$wsReceiving.Range("A14").Value = 12
$wsReceiving.Range("B14").Value = "C.Hollister"
$wsReceiving.Range("C14").Value = 6
$wsReceiving.Range("D14").Value = 4
$wsReceiving.Range("E14").Value = 1
$wsReceiving.Range("F14").Value = 0
$wsReceiving.Range("G14").Value = 2
$wsReceiving.Range("H14").Value = 1

# Row 15 - A.Firkser (shifted down from row 14, stats updated)
$wsReceiving.Range("A15").Value = 13
$wsReceiving.Range("B15").Value = "A.Firkser"
$wsReceiving.Range("C15").Value = 35
$wsReceiving.Range("D15").Value = 30
$wsReceiving.Range("E15").Value = 3
$wsReceiving.Range("F15").Value = 1
$wsReceiving.Range("G15").Value = 5
$wsReceiving.Range("H15").Value = 4

# Row 16 - M.Pruitt (shifted down from row 15, stats updated)
$wsReceiving.Range("A16").Value = 14
$wsReceiving.Range("B16").Value = "M.Pruitt"
$wsReceiving.Range("C16").Value = 14
$wsReceiving.Range("D16").Value = 10
$wsReceiving.Range("E16").Value = 3
$wsReceiving.Range("F16").Value = 2
$wsReceiving.Range("G16").Value = 6
$wsReceiving.Range("H16").Value = 3

# Row 17 - G.Swaim (shifted down from row 16, stats updated)
$wsReceiving.Range("A17").Value = 15
$wsReceiving.Range("B17").Value = "G.Swaim"
$wsReceiving.Range("C17").Value = 33
$wsReceiving.Range("D17").Value = 27
$wsReceiving.Range("E17").Value = 1
$wsReceiving.Range("F17").Value = 1
$wsReceiving.Range("G17").Value = 7
$wsReceiving.Range("H17").Value = 4

# Row 18 - T.Hudson (new row, shifted down from row 17, unchanged stats)
$wsReceiving.Range("A18").Value = 16
$wsReceiving.Range("B18").Value = "T.Hudson"
$wsReceiving.Range("C18").Value = 2
$wsReceiving.Range("D18").Value = 1
$wsReceiving.Range("E18").Value = 1
$wsReceiving.Range("F18").Value = 0
$wsReceiving.Range("G18").Value = 0
$wsReceiving.Range("H18").Value = 0
